$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('NUH NACİ YAZGAN ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=XM1nKVcmbfUtoOteBfvgQPnVlttt06yMfr_cYXKzmRjDSh--IWGUumkwXAa42sRR'),
    @('ONDOKUZ MAYIS ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=dEsunG6lk2T-zEYkWVpizMMPocBPa8fCDfbyxSX9niLGrwTClGJU9GpHP-_2xehJ'),
    @('ORDU ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=705GBkPsPBn5Vaj4UHwca23w0cgf_LIB7WnZbeaefh9bvbgQ2O_T1ouh_7FLbWHJ'),
    @('ORTA DOĞU TEKNİK ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=U89s_CFInFQgY7zTGeGwCdxsIuIquqjF3VvvoAdD6onDSh--IWGUumkwXAa42sRR'),
    @('OSMANİYE KORKUT ATA ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=HgnjgztWgeTfHenW06imS_Lc0RaLkSEhG0BVYp38RYRbvbgQ2O_T1ouh_7FLbWHJ'),
    @('OSTİM TEKNİK ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=5eLCFJ11F4e7oAs72IRQC84Zd294w9NUjaHdBSG23H1eEQelTm-XyRhBW_4oLHrN'),
    @('PAMUKKALE ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=buSE59irWgto5NqYIqKj0GCQaTsNYbpzx8qnY8xDLnTDSh--IWGUumkwXAa42sRR'),
    @('RECEP TAYYİP ERDOĞAN ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=lw4g-cp2QHj3u4IPhXUpMO59JSY3DuCKnYUN-T6QmwzDSh--IWGUumkwXAa42sRR'),
    @('SAKARYA UYGULAMALI BİLİMLER ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=OOZk9nlk7QEJ6URGoB9CphbO-CeDjp0WUsY1guHw7HVZ0rDnHIo0dVL3b_ZK_h53'),
    @('SAKARYA ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=4ssxIM-niFErIBgoRclxqEhzosJXM721rQAK9U9i1WrKRkLi5eSVnKbktSdmDPA1'),
    @('SAMSUN ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=BGTtIj9qMSKPCgQoCwSAVgWrkxPMkWiZaovtn-vWGqrKRkLi5eSVnKbktSdmDPA1'),
    @('SANKO ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=XfwnpyvjRVOV-L160ZO7x1mZYYk7Nf10J_FS6o3bUx0_HK2vOcGPj0bi__hecw25'),
    @('SELÇUK ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=STDouvZ00FlULembVqdSnFLjaL4RoKr-n-47RD0a-fdeEQelTm-XyRhBW_4oLHrN'),
    @('SİİRT ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=ryeBFH608ZJEtzlvnuETVZUc8SqMDGoSCfm7tpzF7AkRd242XOlgWf9_adZszWwP'),
    @('SİNOP ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=zD1wt1ejH1r2p9MnOJ-xdlQ8fbbT-BHaexofjrFh2atbvbgQ2O_T1ouh_7FLbWHJ'),
    @('SİVAS BİLİM VE TEKNOLOJİ ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=2Ul9lq8jZBxrQ-HAjBsj3CPoT6C0AY_NYajpwje8SWu6rGdYDplwTYkFx5fwwQCj'),
    @('SİVAS CUMHURİYET ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=2Ul9lq8jZBxrQ-HAjBsj3KogN0MjqghVNbJxbyBpbQE_HK2vOcGPj0bi__hecw25'),
    @('SÜLEYMAN DEMİREL ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=A4a8ikAWaSJyZ7KpqU8Tu09fS6toC6TViKHPCgybPcZeEQelTm-XyRhBW_4oLHrN'),
    @('ŞIRNAK ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=4upCJx1Uw9pxEWg5WShkGNdNYJWOMS_egkIUxojfGnjGrwTClGJU9GpHP-_2xehJ'),
    @('TARSUS ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=qYXAo1-fPdI9jI-UUrHF0es1p2HgXgFU780rQbNg_7DDSh--IWGUumkwXAa42sRR'),
    @('TED ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=dmTjRIzefxmRIaeNYnE_90OyHKX-xhxVf9Fmqj67n4FeEQelTm-XyRhBW_4oLHrN'),
    @('TEKİRDAĞ NAMIK KEMAL ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=lINch5Jk_HHDgQwVPcA_cC21irjl2G9PABW7Y4ukFblZ0rDnHIo0dVL3b_ZK_h53'),
    @('TOBB EKONOMİ VE TEKNOLOJİ ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=CaTgtXfnKWXajORO0NgHWrseab9EO45CKfDDvMSgU_JeEQelTm-XyRhBW_4oLHrN'),
    @('TOKAT GAZİOSMANPAŞA ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=TV_x5NqUPfAj-baPfTMkupjiqkS70pNNfUCyWGB_Cn1bvbgQ2O_T1ouh_7FLbWHJ'),
    @('TOROS ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=lQIav0Gy9MiZq3-MAGmg70PrpEsM7SNLbNA8FuEy-4PKRkLi5eSVnKbktSdmDPA1'),
    @('TRABZON ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=6Zwrku8aN3f7hEoLyCY8qNR5MQhGepQksim0v3P1z4ZeEQelTm-XyRhBW_4oLHrN'),
    @('TRAKYA ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=FhjgRJNTHexqr_uScLC92ffhIUDsADKJOUAxzzeTPEA_HK2vOcGPj0bi__hecw25'),
    @('TÜRK HAVA KURUMU ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=5ITgxU4tAuKChiflTUqjreamddLYfnAwL30nu_QPT64_HK2vOcGPj0bi__hecw25'),
    @('UFUK ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=beWvCLHYlvqXHPZUbqMXNB20tKZN4m1xQDyAUo-Dw_5Z0rDnHIo0dVL3b_ZK_h53'),
    @('UŞAK ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=beWvCLHYlvqXHPZUbqMXNCLUc7qcsv6BZX-HvbJ15rI_HK2vOcGPj0bi__hecw25'),
    @('VAN YÜZÜNCÜ YIL ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=ZZi-pg0QQb0s8CqAbmRVu5u-mR8UzgJHHZLz_VY-zJXv96FgaKCkCWcpWlro_ywH'),
    @('YALOVA ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=ZAmG36Tv1rjmIerjzTX6F5LngeEsdKY1xLTSmE_TKzxeEQelTm-XyRhBW_4oLHrN'),
    @('YAŞAR ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=CiGTMzhDGOanlv8tpAqEsK_sWyTIlRIfthsBrj0vzZsRd242XOlgWf9_adZszWwP'),
    @('YOZGAT BOZOK ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=CiGTMzhDGOanlv8tpAqEsHqOcM9CVSzjTH62JDxNqxC6rGdYDplwTYkFx5fwwQCj'),
    @('YÜKSEK İHTİSAS ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=1r3Zu71TC9c2XdXIJmxd6J5fwNzuA5A2ffh-gGrpAzFeEQelTm-XyRhBW_4oLHrN'),
    @('ZONGULDAK BÜLENT ECEVİT ÜNİVERSİTESİ', 'https://akademik.yok.gov.tr/AkademikArama/AkademisyenArama?islem=1r3Zu71TC9c2XdXIJmxd6FhBGapc7-c-YXCNg8ohY726rGdYDplwTYkFx5fwwQCj')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

